$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# ---- New shared-string text for row 12 (here-strings preserve newlines/spaces exactly) ----
$sName = @"
Vivek 
"@

$sC = @"
Bhaiya are the most motivating player . Player ke bad 
Senior adviser ki post aapke ke best h mujhe ek bar aapki captaincy me khelna tha vo *sml* me khle liya.😁
Bhaiya bike thodi slow chlaya kro 😅
Day one pe lga ki bhai sab ke kon h in bhaiya ne to bhout attitude h bde gusse vale h  but gradually you are like big brother 🙃
Best wishes for you 
"@

$sD = @"
Aapki captaincy me volleyball iitgn  ne phla tournament jita spoke summit or inter IIT me bhi aapn pre quarters  me pahuche 
Or ha bhaiya mere bhout jyada bolne ka hisab game me pura ho gya tha na  🙃😁
Bhaiya aapka time to time kudasan jana 😁😅

"@

$sE = @"
Bhaiya aapki calmness ka me fan ho gya jab bhi tatu , abhinav bhaiya over aggression me hote the you were the the middle men 
Aapke sath intelectual bate krne me bde mje aate the . Aapki set ki hue ball le taappe marne me bde mje aate the
"@

$sF = @"
Garu bhaiya 🙃😅  aapka har bat pr salah dena .

"@

$sG = @"
Didi jab aapse phli abt mila to me dr dya tha aap is time pr abhinav  bhaiya me gussa kr rhe the   
Then I was like iitne gusse vale senior then uske bad kabhi aapka gussa nhi dekha mene 
Aapka bar bar. Muje ye bolna ki bibek muje hote bhadani h kya kr hu 😆

"@

$sH = @"
 Didi you are like the don 😆😅  aapka vo thada sa guise vala fce usme aap bhout funny lgte ho 🙂

"@

$sI = @"
You are the most cutest senior 🙃😅 
Aapki serve ke bde charche h volleyball samj 
"@

$sJ = @"
 I found you as the most serious senior among all maine aajoo jyada mjak masti krte huee nhi Deka h .
Aapke bina aaki team 
"@

# ---- 1) Capture the current last row (row 11)'s special "bottom border" formatting ----
# (row 11 today uses styles 13/14 - the unique border treatment reserved for the table's last row)
$ws.Range("A11:F11").Copy()

# ---- 2) Grow the table by one row; Excel/ListObject will put this at sheet row 12 ----
$newListRow = $tbl.ListRows.Add()
$newRowIndex = $newListRow.Index + 1

# ---- 3) Paste the captured "last row" formatting onto the new row's populated columns (A-F) ----
$ws.Range("A" + $newRowIndex + ":F" + $newRowIndex).PasteSpecial(-4122)

# Extend the same "last row" cell style (currently on C) across D and G:J as well
$ws.Range("C" + $newRowIndex).Copy()
$ws.Range("D" + $newRowIndex).PasteSpecial(-4122)
$ws.Range("G" + $newRowIndex + ":J" + $newRowIndex).PasteSpecial(-4122)

# Columns K and L have no data in the new row - clear them so no cells are emitted
$ws.Range("K" + $newRowIndex + ":L" + $newRowIndex).Clear()

# ---- 4) Row 11 is no longer the last row, so restore its formatting to the regular striped pattern ----
# (copy from row 9, an existing "odd" data row using the normal, non-last-row border colors)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B9:C9").Copy()
$ws.Range("B11:C11").PasteSpecial(-4122)
$ws.Range("E9:F9").Copy()
$ws.Range("E11:F11").PasteSpecial(-4122)

# ---- 5) Populate the new row's values ----
$ws.Cells.Item($newRowIndex, 1).Value = 45767.84058173611
$ws.Cells.Item($newRowIndex, 2).Value = $sName
$ws.Cells.Item($newRowIndex, 3).Value = $sC
$ws.Cells.Item($newRowIndex, 4).Value = $sD
$ws.Cells.Item($newRowIndex, 5).Value = $sE
$ws.Cells.Item($newRowIndex, 6).Value = $sF
$ws.Cells.Item($newRowIndex, 7).Value = $sG
$ws.Cells.Item($newRowIndex, 8).Value = $sH
$ws.Cells.Item($newRowIndex, 9).Value = $sI
$ws.Cells.Item($newRowIndex, 10).Value = $sJ

# Multi-line text auto-expands the row height; restore it to the sheet's standard height
# so no explicit ht/customHeight attribute is written (matching the other rows).
$ws.Rows.Item($newRowIndex).AutoFit()
